# This workbook stores a daily price log for "Espinaca" (Hortaliza) at the
# Vega Central Mapocho de Santiago market. A new day's record is inserted
# at row 553, which pushes the existing rows 553-679 down to 554-680
# (dimension grows from A1:R679 to A1:R680).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 553; Excel shifts rows 553:679 down to 554:680 and
# copies the row formatting (including the date-cell number format) from
# the row above, exactly like a native "Insert Row" in the UI.
$ws.Rows.Item(553).Insert()

# Populate the newly inserted row 553 with the new daily record.
$ws.Cells.Item(553, 1).Value = 9
$ws.Cells.Item(553, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(553, 3).Value = "Metropolitana"
$ws.Cells.Item(553, 4).Value = 45211
$ws.Cells.Item(553, 5).Value = 13
$ws.Cells.Item(553, 6).Value = 100112012
$ws.Cells.Item(553, 7).Value = "Espinaca"
$ws.Cells.Item(553, 8).Value = "Sin especificar"
$ws.Cells.Item(553, 9).Value = "Primera"
$ws.Cells.Item(553, 10).Value = 160
$ws.Cells.Item(553, 11).Value = 7000
$ws.Cells.Item(553, 12).Value = 8000
$ws.Cells.Item(553, 13).Value = 7500
$ws.Cells.Item(553, 14).Value = "$/cuna 10 kilos"
$ws.Cells.Item(553, 15).Value = "Provincia de Chacabuco"
$ws.Cells.Item(553, 16).Value = 750
$ws.Cells.Item(553, 17).Value = 10
$ws.Cells.Item(553, 18).Value = "Hortaliza"

# Make sure the date cell keeps the same date/time number format used by
# the rest of the column (Insert() should already propagate it, but set
# it explicitly to be safe).
$ws.Cells.Item(553, 4).NumberFormat = $ws.Cells.Item(554, 4).NumberFormat()
